$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had 7 data rows (1 header + 6 data rows): rows 2-4 had "ECs" as
# the sending cluster and rows 5-7 had "FAPs" as the sending cluster (both against the
# same Dll3/Notch4 ligand-receptor pair, one row per target cluster: ECs, FAPs, MuSCs).
# With the updated TPM numbers, the "ECs" sending-cluster rows are dropped entirely, so
# delete rows 5-7 and turn what used to be rows 2-4 ("ECs" sender) into the "FAPs" sender
# rows with the newly recalculated values below.
$ws.Rows("5:7").Delete()

# Row 2: update cells whose value changed with the new TPM-based recalculation
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("G2").Value2 = 0.2138853333333333
$ws.Range("H2").Value2 = 0.641656
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("M2").Value2 = 31.618405
$ws.Range("N2").Value2 = 94.855215
$ws.Range("O2").Value2 = 0.8578613706944929
$ws.Range("P2").Value2 = 0.8578613706944929
$ws.Range("Q2").Value2 = 6.762713092893334
$ws.Range("R2").Value2 = 60.86441783604
$ws.Range("S2").Value2 = 0.8578613706944929
$ws.Range("T2").Value2 = 0.8578613706944929

# Row 3: update cells whose value changed with the new TPM-based recalculation
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("G3").Value2 = 0.2138853333333333
$ws.Range("H3").Value2 = 0.641656
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = 1
$ws.Range("O3").Value2 = 0.08747555172986397
$ws.Range("P3").Value2 = 0.08747555172986396
$ws.Range("Q3").Value2 = 0.6895893429875555
$ws.Range("R3").Value2 = 6.206304086887999
$ws.Range("S3").Value2 = 0.08747555172986397
$ws.Range("T3").Value2 = 0.08747555172986396

# Row 4: update cells whose value changed with the new TPM-based recalculation
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("G4").Value2 = 0.2138853333333333
$ws.Range("H4").Value2 = 0.641656
$ws.Range("I4").Value2 = 1
$ws.Range("J4").Value2 = 1
$ws.Range("M4").Value2 = 2.014730333333334
$ws.Range("N4").Value2 = 6.044191000000001
$ws.Range("O4").Value2 = 0.05466307757564324
$ws.Range("P4").Value2 = 0.05466307757564324
$ws.Range("Q4").Value2 = 0.4309212689217778
$ws.Range("R4").Value2 = 3.878291420296
$ws.Range("S4").Value2 = 0.05466307757564324
$ws.Range("T4").Value2 = 0.05466307757564324
